$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 33: label "doWork03" already in A33; add a new column header in B33.
$ws.Range("B33").Value = "Buffered Input Stream"

# Rows 34-43 used to hold a single column of text labels (A34:A43).
# Replace them with a full numeric data table spanning A:H, matching the
# layout used by the other statistics blocks on the sheet.
$data = @(
  @(94.4,                34.442,               72.9,                11.295,              70.1,                27.286,              38.3,                10.292),
  @(65.7,                2.0099999999999998,   40.799999999999997,  10.275,              58.7,                17.66,               34.7,                0.53800000000000003),
  @(66.400000000000006,  2.4529999999999998,   36.5,                0.47,                42.7,                2.101,               34.200000000000003, 0.69),
  @(65.900000000000006,  2.1720000000000002,   36.6,                1.7509999999999999, 47,                  7.7279999999999998, 34.5,                0.61599999999999999),
  @(65.3,                2.375,                36.6,                1.53,                46.7,                13.627000000000001, 34.700000000000003, 0.73),
  @(65.3,                2.101,                36.700000000000003,  1.65,                40.799999999999997,  2.056,               34.6,                0.60699999999999998),
  @(67.3,                3.8260000000000001,   36,                  1.1819999999999999, 39.799999999999997,  1.6120000000000001, 34.799999999999997, 0.49299999999999999),
  @(65.099999999999994,  1.357,                36.200000000000003,  1.3340000000000001, 39.9,                1.2,                 34.700000000000003, 0.59799999999999998),
  @(65.8,                1.7150000000000001,   36.1,                0.95799999999999996, 39.9,               1.6850000000000001, 34.6,                0.56799999999999995),
  @(65.599999999999994,  2.218,                36.4,                1.359,               38.700000000000003, 2.8029999999999999, 34.799999999999997, 0.70299999999999996)
)

$startRow = 34
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  for ($c = 0; $c -lt $data[$i].Count; $c++) {
    $col = $c + 1
    $ws.Cells.Item($row, $col).Value = $data[$i][$c]
  }
}

# Update the view: current selection moved to the newly filled block.
$ws.Activate()
$ws.Range("G34:G43").Select()
